$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "Bug ID" number-cell formatting (bold font + border +
# centered/top alignment, style index 1 in the original sheet) down onto the
# three new rows so A11:A13 keep the same look as A2:A10.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A13").PasteSpecial(-4122) | Out-Null

# Row 11 - "ffg not working" (still open, no update/close date yet)
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "ffg not working"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = "Open"
$ws.Range("E11").Value = "2023-07-28 12:38:00"
$ws.Range("F11").Value = ""

# Row 12 - "io none" (still open, no update/close date yet)
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "io none"
$ws.Range("C12").Value = "io none"
$ws.Range("D12").Value = "Open"
$ws.Range("E12").Value = "2023-07-28 12:41:13"
$ws.Range("F12").Value = ""

# Row 13 - "er fgv" (fixed + passed)
$ws.Range("A13").Value = 13
$ws.Range("B13").Value = "er fgv"
$ws.Range("C13").Value = "er fgv fixed"
$ws.Range("D13").Value = "`n            Passed"
$ws.Range("E13").Value = "2023-07-28 16:13:23"
$ws.Range("F13").Value = "2023-07-28 16:15:07"
